# Fruta / hortaliza, semanal
#
# A new weekly price record (Fecha 2021-11-22, serial 44522) is inserted
# for "Feria Lagunitas de Puerto Montt - Apio" at row 136, right before the
# existing row for Fecha 2021-02-19 (serial 44246). All the rows that used
# to be at 136-162 shift down by one to 137-163, and the sheet's used range
# grows from A1:R162 to A1:R163.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh blank row at 136, pushing the old row 136 (and everything
# below it) down to row 137, etc.
$ws.Rows("136:136").Insert()

# The record that used to live at row 137 is now (after the shift) at row
# 138. Duplicate it back into the newly-created row 136 so every column
# (including styles) matches, then overwrite just the date to the new
# week's value.
$ws.Range("A138:R138").Copy()
$ws.Range("A136:R136").PasteSpecial(-4104)
$ws.Range("D136").Value = 44522
